$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7107720971107483
$ws.Range("B1").Value = 1.459053039550781
$ws.Range("C1").Value = 4.033444404602051
$ws.Range("D1").Value = 2.704841136932373
$ws.Range("E1").Value = 0.5901511311531067
